$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily dates appended by the "MV -datos-" update.
$dates = @("29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021", "03-11-2021")

$startRow = 303
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Write the date as literal text. Some of the new dates (e.g. 01-11-2021)
    # look like valid dates to Excel's locale-aware parser and would silently
    # be converted to a date serial if assigned directly via .Value. Building
    # it through a formula that returns a text string, then collapsing the
    # formula down to its static value with a values-only paste, keeps the
    # cell as plain text (t="s") without touching any cell formatting/styles.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Formula = '="' + $dates[$i] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = 12836
    $ws.Cells.Item($row, 3).Value = 0

    # The very last appended row (03-11-2021) has not yet received a value
    # for column D in the source data.
    if ($row -lt ($startRow + $dates.Length - 1)) {
        $ws.Cells.Item($row, 4).Value = 393
    }
}

$excel.Application.CutCopyMode = $false
